$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.245.95'
$ws.Range('E2').Value = '  -3.64%  '
$ws.Range('D3').Value = '3.143.14'
$ws.Range('E3').Value = '  -3.11%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = '''603.12'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.30%  '
$ws.Range('D6').Value = '''146.18'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -7.05%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '3.144.09'
$ws.Range('E8').Value = '  -3.00%  '
$ws.Range('E9').Value = '  -4.31%  '
$ws.Range('E10').Value = '  -7.53%  '
$ws.Range('D11').Value = '''5.49'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -5.44%  '
$ws.Range('E12').Value = '  -5.90%  '
$ws.Range('D13').Value = '''0.0000250'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -8.17%  '
$ws.Range('D14').Value = '''35.88'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -8.27%  '
$ws.Range('D15').Value = '3.657.67'
$ws.Range('E15').Value = '  -3.22%  '
$ws.Range('D16').Value = '64.250.97'
$ws.Range('E16').Value = '  -3.67%  '
$ws.Range('E17').Value = '  +0.39%  '
$ws.Range('D18').Value = '3.142.15'
$ws.Range('E18').Value = '  -3.34%  '
$ws.Range('D19').Value = '''6.91'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.44%  '
$ws.Range('D20').Value = '''477.09'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -6.21%  '
$ws.Range('D21').Value = '''14.59'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.96%  '
$ws.Range('D22').Value = '''0.706'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.28%  '
$ws.Range('D23').Value = '''7.67'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.87%  '
$ws.Range('D24').Value = '''13.64'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -7.01%  '
$ws.Range('D25').Value = '''83.24'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.39%  '
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('E27').Value = '  -5.21%  '
$ws.Range('D28').Value = '''8.36'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -7.94%  '
$ws.Range('D29').Value = '''2.18'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -7.42%  '
$ws.Range('D30').Value = '''6.71'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.75%  '
$ws.Range('D31').Value = '''0.112'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -36.99%  '
$ws.Range('E32').Value = '  +0.04%  '
$ws.Range('E33').Value = '  -5.79%  '
$ws.Range('D34').Value = '''26.07'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -7.85%  '
$ws.Range('E35').Value = '  -5.09%  '
$ws.Range('D36').Value = '''54.15'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.07%  '
$ws.Range('E37').Value = '  -6.50%  '
$ws.Range('D38').Value = '0.0₃0709'
$ws.Range('E38').Value = '  -11.89%  '
$ws.Range('D39').Value = '''442.51'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -10.40%  '
$ws.Range('E40').Value = '  -12.51%  '
$ws.Range('D41').Value = '''0.0394'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -7.47%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').Value = '''0.118'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -7.91%  '
$ws.Range('B43').Value = 'Cosmos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D43').Value = '''8.40'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.00%  '
$ws.Range('D44').Value = '2.821.31'
$ws.Range('E44').Value = '  -4.22%  '
$ws.Range('E45').Value = '  -9.96%  '
$ws.Range('E46').Value = '  -8.96%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value = '''26.27'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -7.25%  '
$ws.Range('B48').Value = 'USDe'
$ws.Range('C48').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D48').Value = '''0.998'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.14%  '
$ws.Range('E49').Value = '  -4.68%  '
$ws.Range('D50').Value = '''2.29'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.33%  '
$ws.Range('D51').Value = '''117.53'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.22%  '
